# "Generate Report for Handback" - update the localization-status workbook
# to reflect that both the zh-cn and de-de handback packages have been
# generated: the per-language sheets get their target/handback file
# hyperlinks + handback datetime filled in, the Overview sheet's status
# columns flip from "Ready for handoff" to "Handed back: in sync with
# en-US", and a couple of columns are widened so the new longer file-name
# values are readable.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ba0694461c6fce7f86e920b4ec6cc9368bdd6519/e2e/"

$file1Name = "24497038-753b-48ff-961d-054ed183cbdc.md"
$file2Name = "84bfd5b9-4cc0-4434-a83d-07728a2be18e.md"

$file1Url = $ghBase + $file1Name
$file2Url = $ghBase + $file2Name

# ---------------------------------------------------------------------
# Overview sheet: flip the per-language status cells to the new text and
# widen the zh-cn / de-de columns so the longer status text fits.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): fill in the "Latest Target File"
# (I) and "Latest Handback File" (J) hyperlinked file names, stamp the
# "Latest Handback DateTime" (K), refresh the Status column (C) text,
# and widen columns C/I/J to fit the new content.
# ---------------------------------------------------------------------
$languages = @(
    @{ Name = "zh-cn"; Suffix = "zh-cn.xlf"; HandbackTime = "2016-10-19 17:29:23" },
    @{ Name = "de-de"; Suffix = "de-de.xlf"; HandbackTime = "2016-10-19 17:29:41" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Name)

    $xlf1 = "24497038-753b-48ff-961d-054ed183cbdc.5d0b70896959f4334ab5824ed556404a867a235e." + $lang.Suffix
    $xlf2 = "84bfd5b9-4cc0-4434-a83d-07728a2be18e.ea426d2b0e5ccb5c88ef19b34ef354159711b5a1." + $lang.Suffix

    # Status column refresh
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Handback datetime stamps
    $ws.Range("K2").Value = $lang.HandbackTime
    $ws.Range("K3").Value = $lang.HandbackTime

    # Latest Handback File (xliff) names
    $ws.Range("J2").Value = $xlf1
    $ws.Range("J3").Value = $xlf2

    # Rebuild the hyperlinks: row 2's target-file link, then row 3's
    # source-file link followed by row 3's target-file link, so the
    # relationship ids land in the same order as the source file link
    # that already exists for row 2.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $file1Url, "", "", $file1Name)
    $ws.Hyperlinks.Add($ws.Range("I2"), $file1Url, "", "", $file1Name)
    $ws.Hyperlinks.Add($ws.Range("A3"), $file2Url, "", "", $file2Name)
    $ws.Hyperlinks.Add($ws.Range("I3"), $file2Url, "", "", $file2Name)

    # Column widths for the newly-populated / widened columns.
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}
